$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column H values (shared strings, in this order: "Story id", "MOSIP-9800", "MOSIP-11479", "MOSIP-11986")
$ws.Cells.Item(1, 8).Value = "Story id"
$ws.Cells.Item(2, 8).Value = "MOSIP-9800"
$ws.Cells.Item(3, 8).Value = "MOSIP-9800"
$ws.Cells.Item(4, 8).Value = "MOSIP-11479"
$ws.Cells.Item(5, 8).Value = "MOSIP-11479"
$ws.Cells.Item(6, 8).Value = "MOSIP-11479"
$ws.Cells.Item(7, 8).Value = "MOSIP-11479"
$ws.Cells.Item(8, 8).Value = "MOSIP-11479"
$ws.Cells.Item(9, 8).Value = "MOSIP-11479"
$ws.Cells.Item(10, 8).Value = "MOSIP-11986"

# Header cell H1 - match the header formatting used by the rest of row 1 (copy format from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column width for the new column
$ws.Columns.Item(8).ColumnWidth = 16.5

# Move the active selection (no longer on the last row)
$ws.Range("J5").Select()
